$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 48, shifting the existing rows 48:85 down to 49:86.
$ws.Rows("48:48").Insert()

# Populate the newly inserted row 48 with the new weekly price record.
$ws.Range("A48").Value = 9
$ws.Range("B48").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C48").Value = "Metropolitana"
$ws.Range("D48").Value = 44606
$ws.Range("E48").Value = 13
$ws.Range("F48").Value = 100114007
$ws.Range("G48").Value = "Jengibre"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 520
$ws.Range("K48").Value = 17000
$ws.Range("L48").Value = 18000
$ws.Range("M48").Value = 17500
$ws.Range("N48").Value = "$/caja 13 kilos"
$ws.Range("O48").Value = "Perú"
$ws.Range("P48").Value = 1346
$ws.Range("Q48").Value = 13
$ws.Range("R48").Value = "Hortaliza"
